$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Target values for rows 34-40, columns A,B,D,E,F,G,H,Q,R
# This corresponds to a cyclic reshuffle of species records among rows 34,36,39 and 35,38,40 (row 37 unchanged)
$data = @{
    34 = @{ "A" = 112038596; "B" = 90087; "D" = 'LC'; "E" = 3298; "F" = 'Trådticka'; "G" = 'Climacocystis borealis'; "H" = '(Fr.) Kotl. & Pouzar'; "Q" = 616076.0611235843; "R" = 6895427.595461337 }
    35 = @{ "A" = 112038603; "B" = 89369; "D" = 'LC'; "E" = 5447; "F" = 'Vedticka'; "G" = 'Fuscoporia viticola'; "H" = '(Schwein.) Murrill'; "Q" = 615968.1934313668; "R" = 6895405.650930508 }
    36 = @{ "A" = 112038599; "B" = 89423; "D" = 'NT'; "E" = 5432; "F" = 'Granticka'; "G" = 'Porodaedalea chrysoloma'; "H" = '(Fr.) Fiasson & Niemelä'; "Q" = 616070.2961488151; "R" = 6895499.860901954 }
    37 = @{ "A" = 112038600; "B" = 86223; "D" = 'NT'; "E" = 4412; "F" = 'Äggvaxskivling'; "G" = 'Hygrophorus karstenii'; "H" = 'Sacc. & Cub.'; "Q" = 616034.1211971109; "R" = 6895585.10294092 }
    38 = @{ "A" = 112038602; "B" = 86223; "D" = 'NT'; "E" = 4412; "F" = 'Äggvaxskivling'; "G" = 'Hygrophorus karstenii'; "H" = 'Sacc. & Cub.'; "Q" = 616026.2967975155; "R" = 6895553.979090866 }
    39 = @{ "A" = 112038601; "B" = 73634; "D" = 'LC'; "E" = 6426; "F" = 'Kattfotslav'; "G" = 'Felipes leucopellaeus'; "H" = '(Ach.) Frisch & G.Thor'; "Q" = 616012.5978259755; "R" = 6895611.944218947 }
    40 = @{ "A" = 112038604; "B" = 89845; "D" = 'VU'; "E" = 1209; "F" = 'Rynkskinn'; "G" = 'Phlebia centrifuga'; "H" = 'P.Karst.'; "Q" = 615977.7276359925; "R" = 6895550.438170813 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
